# Applies the resume content edits described by the commit diff.
$d = $word.ActiveDocument

# 1. Professional summary: replace "producing reports and helping users
#    understand issues." with "team management and technical aspects."
#    (this also naturally drops the gramStart/gramEnd proofErr wrapping
#    around the old "reports" run, since that text no longer exists).
$d.Content.Find.Execute(
    "producing reports and helping users understand issues.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "team management and technical aspects.", 2) | Out-Null

# 2. Skills table: "Telerik Rad Controls" -> "Angular"
$d.Content.Find.Execute(
    "Telerik Rad Controls",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Angular", 2) | Out-Null

# 3. "...that includes: Design," -> "...that includes Design," (drop colon)
$d.Content.Find.Execute(
    "that includes: Design,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "that includes Design,", 2) | Out-Null

# 4. Remove the whole bullet paragraph "Championed and implemented best
#    practices and delivery standards."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Championed and implemented best practices and delivery standards.*") {
        $p.Range.Delete()
        break
    }
}

# 5. "It is truly develop using the WEB API." -> "It is develop using the WEB API."
$d.Content.Find.Execute(
    "It is truly develop using the WEB API.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "It is develop using the WEB API.", 2) | Out-Null

# 6. "...service it is act as a Gateway service." -> "...service and act as an Gateway service."
$d.Content.Find.Execute(
    "service it is act as a Gateway service.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "service and act as an Gateway service.", 2) | Out-Null

# 7. Drop the trailing "This will get all the data from the database and
#    send the response to the" tail sentence (including the "response"
#    badword run), leaving the paragraph ending at "...respectively ".
$d.Content.Find.Execute(
    "respectively. This will get all the data from the database and send the response to the ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "respectively ", 2) | Out-Null

# 8. "This is and UI" -> "This is an UI" (typo fix)
$d.Content.Find.Execute(
    "This is and UI",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This is an UI", 2) | Out-Null

# 9. "...align those requests in proper way. Can match..." ->
#    "...align those requests in proper way with the respective view model
#    using MVVM pattern. Can match..."
$d.Content.Find.Execute(
    "align those requests in proper way. Can match",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "align those requests in proper way with the respective view model using MVVM pattern. Can match", 2) | Out-Null
